$wb = $excel.ActiveWorkbook

$survey = $wb.Worksheets.Item("survey")
$choices = $wb.Worksheets.Item("choices")

# --- "survey" sheet: insert a new 3-row "agriculture" external-link block ---
# immediately before the existing first block ("geotagger", currently rows 8:10),
# matching the XLSForm "branch_label / external_link / exit section" pattern
# used by every other test form in this sheet.
$survey.Rows("8:10").Insert()

# Clone the formatting of the (now shifted) geotagger block into the new blank rows.
$survey.Range("A11").Copy()
$survey.Range("A8").PasteSpecial(-4122)   # xlPasteFormats
$survey.Range("B11").Copy()
$survey.Range("B8").PasteSpecial(-4122)

$survey.Range("A12").Copy()
$survey.Range("A9").PasteSpecial(-4122)
$survey.Range("B12").Copy()
$survey.Range("B9").PasteSpecial(-4122)

$survey.Range("A13").Copy()
$survey.Range("A10").PasteSpecial(-4122)
$survey.Range("B13").Copy()
$survey.Range("B10").PasteSpecial(-4122)

$survey.Rows(8).RowHeight = 17
$survey.Rows(9).RowHeight = 59
$survey.Rows(10).RowHeight = 17

$survey.Range("A8").Value = "agriculture"
$survey.Range("B9").Value = "'?' + opendatakit.getHashString('../tables/agriculture/forms/agriculture/',null)"
$survey.Range("E9").Value = "external_link"
$survey.Range("G9").Value = "Open form"
$survey.Range("C10").Value = "exit section"

# --- "choices" sheet: add the matching "agriculture" choice row ---
# right before the existing last row ("geoweather", currently row 20).
$choices.Rows(20).Insert()
$choices.Range("A20").Value = "test_forms"
$choices.Range("B20").Value = "agriculture"
$choices.Range("C20").Value = "Agriculture"

# --- restore selections / active sheet to match the saved workbook state ---
$survey.Range("A9").Select()
$choices.Range("C25").Select()
$choices.Activate()

Write-Output "done"
